$wb = $excel.ActiveWorkbook
try {
  Write-Host "WindowState:" $excel.ActiveWindow.WindowState
} catch { Write-Host "err windowstate:" $_ }
try {
  Write-Host "Width:" $excel.ActiveWindow.Width
  Write-Host "Height:" $excel.ActiveWindow.Height
} catch { Write-Host "err width/height:" $_ }
